# Design_Table_T2.docx — add "with Bonferroni adjustment for multiple testing"
# after every "... pairwise contrasts are calculated using pairs()" sentence,
# as a tracked insertion authored by Alex Strobel.

$word.UserName = "Alex Strobel"

$d = $word.ActiveDocument
$d.TrackRevisions = $true

$insertText = " with Bonferroni adjustment for multiple testing"
# The 5th occurrence of the sentence straddles a page break across two runs;
# there the word order comes out as "...pairs()" + "with ... testing " + "."
# (no space before "with", a trailing space before the final period).
$insertTextNoLeadingSpace = "with Bonferroni adjustment for multiple testing "

$search = $d.Content
$search.Start = 0
$occurrence = 0

while ($search.Find.Execute("pairs().", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $occurrence = $occurrence + 1

    # $search now spans "pairs()." (period included). Drop the trailing period
    # from the match so the insertion point sits right before it.
    $search.MoveEnd(1, -1) | Out-Null
    $insertAt = $search.End

    $point = $d.Range($insertAt, $insertAt)
    if ($occurrence -eq 5) {
        $point.InsertAfter($insertTextNoLeadingSpace)
    } else {
        $point.InsertAfter($insertText)
    }

    # Resume searching right after the (still untouched) trailing period so we
    # don't re-match the text we just inserted and don't skip anything.
    $resumeAt = $search.End + 1
    $search.Start = $resumeAt
    $search.End = $resumeAt
}

Write-Output "Bonferroni sentence updated in $occurrence place(s)."
